# Generate Report for Handoff
# The localization run finished: the handoff package's xliff/markdown file
# GUID rolled over from e2b6ad54-... to 7a76190d-..., and the handoff/
# generate timestamps advanced a few seconds. Propagate the new file names
# and timestamps into the report workbook (Overview + per-locale sheets),
# including the hyperlinks that point at the handed-off markdown file.

$wb = $excel.ActiveWorkbook

$newGuid = "7a76190d-de93-40ea-8880-55568815d466"

$newZhXlf = $newGuid + ".09e8e78fff570da572b61e0f464184ed810dbd09.zh-cn.xlf"
$newDeXlf = $newGuid + ".09e8e78fff570da572b61e0f464184ed810dbd09.de-de.xlf"

$newMdName = $newGuid + ".md"
$newMdPath = "e2e\" + $newGuid + ".md"

$newGenerateDate = "2016-08-31 07:06:31"
$newZhHandoffDate = "2016-08-31 07:06:27"

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7a638f64e58398f4bc90a0d06017a5ce8448025/e2e/"
$newGithubUrl = $githubBase + $newMdName

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newGenerateDate

# Keep the B2 hyperlink's visible text (and target) in sync with the renamed file.
$hlOverview = $wsOverview.Hyperlinks.Item(1)
$hlOverview.Address = $newGithubUrl
$hlOverview.TextToDisplay = $newMdPath

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate

$hlZh = $wsZh.Hyperlinks.Item(1)
$hlZh.Address = $newGithubUrl
$hlZh.TextToDisplay = $newMdName

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newGenerateDate

$hlDe = $wsDe.Hyperlinks.Item(1)
$hlDe.Address = $newGithubUrl
$hlDe.TextToDisplay = $newMdName
